$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5725316666666668
$ws.Range("H2").Value = 1.717595
$ws.Range("I2").Value = 0.3864899584549088
$ws.Range("J2").Value = 0.3864899584549088
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 14.861848
$ws.Range("N2").Value = 44.585544
$ws.Range("O2").Value = 0.09055189482833943
$ws.Range("P2").Value = 0.09055189482833945
$ws.Range("Q2").Value = 8.508878605186668
$ws.Range("R2").Value = 76.57990744668001
$ws.Range("S2").Value = 0.03499739807021818
$ws.Range("T2").Value = 0.03499739807021818

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5725316666666668
$ws.Range("H3").Value = 1.717595
$ws.Range("I3").Value = 0.3864899584549088
$ws.Range("J3").Value = 0.3864899584549088
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 31.07813833333333
$ws.Range("N3").Value = 93.234415
$ws.Range("O3").Value = 0.1893562842131466
$ws.Range("P3").Value = 0.1893562842131466
$ws.Range("Q3").Value = 17.79321833688056
$ws.Range("R3").Value = 160.138965031925
$ws.Range("S3").Value = 0.07318430241871492
$ws.Range("T3").Value = 0.07318430241871493

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5725316666666668
$ws.Range("H4").Value = 1.717595
$ws.Range("I4").Value = 0.3864899584549088
$ws.Range("J4").Value = 0.3864899584549088
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 18.10188466666667
$ws.Range("N4").Value = 54.305654
$ws.Range("O4").Value = 0.1102931450066459
$ws.Range("P4").Value = 0.1102931450066459
$ws.Range("Q4").Value = 10.36390219801445
$ws.Range("R4").Value = 93.27511978213002
$ws.Range("S4").Value = 0.0426271930314798
$ws.Range("T4").Value = 0.04262719303147981

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5725316666666668
$ws.Range("H5").Value = 1.717595
$ws.Range("I5").Value = 0.3864899584549088
$ws.Range("J5").Value = 0.3864899584549088
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 100.0833306666667
$ws.Range("N5").Value = 300.249992
$ws.Range("O5").Value = 0.609798675951868
$ws.Range("P5").Value = 0.6097986759518681
$ws.Range("Q5").Value = 57.30087611213778
$ws.Range("R5").Value = 515.70788500924
$ws.Range("S5").Value = 0.2356810649344959
$ws.Range("T5").Value = 0.2356810649344959

$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.3064436666666667
$ws.Range("H6").Value = 0.919331
$ws.Range("I6").Value = 0.2068661122070742
$ws.Range("J6").Value = 0.2068661122070743
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 14.861848
$ws.Range("N6").Value = 44.585544
$ws.Range("O6").Value = 0.09055189482833943
$ws.Range("P6").Value = 0.09055189482833945
$ws.Range("Q6").Value = 4.554319194562667
$ws.Range("R6").Value = 40.988872751064
$ws.Range("S6").Value = 0.01873211843612245
$ws.Range("T6").Value = 0.01873211843612246

$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.3064436666666667
$ws.Range("H7").Value = 0.919331
$ws.Range("I7").Value = 0.2068661122070742
$ws.Range("J7").Value = 0.2068661122070743
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 31.07813833333333
$ws.Range("N7").Value = 93.234415
$ws.Range("O7").Value = 0.1893562842131466
$ws.Range("P7").Value = 0.1893562842131466
$ws.Range("Q7").Value = 9.523698664040555
$ws.Range("R7").Value = 85.713287976365
$ws.Range("S7").Value = 0.03917139833715142
$ws.Range("T7").Value = 0.03917139833715143

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.3064436666666667
$ws.Range("H8").Value = 0.919331
$ws.Range("I8").Value = 0.2068661122070742
$ws.Range("J8").Value = 0.2068661122070743
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 18.10188466666667
$ws.Range("N8").Value = 54.305654
$ws.Range("O8").Value = 0.1102931450066459
$ws.Range("P8").Value = 0.1102931450066459
$ws.Range("Q8").Value = 5.547207910830444
$ws.Range("R8").Value = 49.924871197474
$ws.Range("S8").Value = 0.02281591411061592
$ws.Range("T8").Value = 0.02281591411061593

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.3064436666666667
$ws.Range("H9").Value = 0.919331
$ws.Range("I9").Value = 0.2068661122070742
$ws.Range("J9").Value = 0.2068661122070743
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 100.0833306666667
$ws.Range("N9").Value = 300.249992
$ws.Range("O9").Value = 0.609798675951868
$ws.Range("P9").Value = 0.6097986759518681
$ws.Range("Q9").Value = 30.66990282170578
$ws.Range("R9").Value = 276.029125395352
$ws.Range("S9").Value = 0.1261466813231844
$ws.Range("T9").Value = 0.1261466813231845

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5451493333333333
$ws.Range("H10").Value = 1.635448
$ws.Range("I10").Value = 0.3680053968340403
$ws.Range("J10").Value = 0.3680053968340404
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 14.861848
$ws.Range("N10").Value = 44.585544
$ws.Range("O10").Value = 0.09055189482833943
$ws.Range("P10").Value = 0.09055189482833945
$ws.Range("Q10").Value = 8.101926529301332
$ws.Range("R10").Value = 72.91733876371198
$ws.Range("S10").Value = 0.03332358599037734
$ws.Range("T10").Value = 0.03332358599037735

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5451493333333333
$ws.Range("H11").Value = 1.635448
$ws.Range("I11").Value = 0.3680053968340403
$ws.Range("J11").Value = 0.3680053968340404
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 31.07813833333333
$ws.Range("N11").Value = 93.234415
$ws.Range("O11").Value = 0.1893562842131466
$ws.Range("P11").Value = 0.1893562842131466
$ws.Range("Q11").Value = 16.94222639365778
$ws.Range("R11").Value = 152.48003754292
$ws.Range("S11").Value = 0.06968413451487833
$ws.Range("T11").Value = 0.06968413451487834

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5451493333333333
$ws.Range("H12").Value = 1.635448
$ws.Range("I12").Value = 0.3680053968340403
$ws.Range("J12").Value = 0.3680053968340404
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 18.10188466666667
$ws.Range("N12").Value = 54.305654
$ws.Range("O12").Value = 0.1102931450066459
$ws.Range("P12").Value = 0.1102931450066459
$ws.Range("Q12").Value = 9.868230358110221
$ws.Range("R12").Value = 88.81407322299199
$ws.Range("S12").Value = 0.04058847259624507
$ws.Range("T12").Value = 0.04058847259624509

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5451493333333333
$ws.Range("H13").Value = 1.635448
$ws.Range("I13").Value = 0.3680053968340403
$ws.Range("J13").Value = 0.3680053968340404
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 100.0833306666667
$ws.Range("N13").Value = 300.249992
$ws.Range("O13").Value = 0.609798675951868
$ws.Range("P13").Value = 0.6097986759518681
$ws.Range("Q13").Value = 54.56036099071287
$ws.Range("R13").Value = 491.0432489164159
$ws.Range("S13").Value = 0.2244092037325395
$ws.Range("T13").Value = 0.2244092037325396

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.05723766666666667
$ws.Range("H14").Value = 0.171713
$ws.Range("I14").Value = 0.03863853250397663
$ws.Range("J14").Value = 0.03863853250397663
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 14.861848
$ws.Range("N14").Value = 44.585544
$ws.Range("O14").Value = 0.09055189482833943
$ws.Range("P14").Value = 0.09055189482833945
$ws.Range("Q14").Value = 0.8506575018746667
$ws.Range("R14").Value = 7.655917516872
$ws.Range("S14").Value = 0.003498792331621466
$ws.Range("T14").Value = 0.003498792331621467

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.05723766666666667
$ws.Range("H15").Value = 0.171713
$ws.Range("I15").Value = 0.03863853250397663
$ws.Range("J15").Value = 0.03863853250397663
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 31.07813833333333
$ws.Range("N15").Value = 93.234415
$ws.Range("O15").Value = 0.1893562842131466
$ws.Range("P15").Value = 0.1893562842131466
$ws.Range("Q15").Value = 1.778840122543889
$ws.Range("R15").Value = 16.009561102895
$ws.Range("S15").Value = 0.0073164489424019
$ws.Range("T15").Value = 0.007316448942401903

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.05723766666666667
$ws.Range("H16").Value = 0.171713
$ws.Range("I16").Value = 0.03863853250397663
$ws.Range("J16").Value = 0.03863853250397663
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 18.10188466666667
$ws.Range("N16").Value = 54.305654
$ws.Range("O16").Value = 0.1102931450066459
$ws.Range("P16").Value = 0.1102931450066459
$ws.Range("Q16").Value = 1.036109640589111
$ws.Range("R16").Value = 9.324986765302
$ws.Range("S16").Value = 0.004261565268305094
$ws.Range("T16").Value = 0.004261565268305096

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.05723766666666667
$ws.Range("H17").Value = 0.171713
$ws.Range("I17").Value = 0.03863853250397663
$ws.Range("J17").Value = 0.03863853250397663
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 100.0833306666667
$ws.Range("N17").Value = 300.249992
$ws.Range("O17").Value = 0.609798675951868
$ws.Range("P17").Value = 0.6097986759518681
$ws.Range("Q17").Value = 5.728536319588444
$ws.Range("R17").Value = 51.55682687629599
$ws.Range("S17").Value = 0.02356172596164816
$ws.Range("T17").Value = 0.02356172596164817

